$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" banner text ---
$ws.Range("B3").Value = "Updated on 3 july 2021"

# --- Add new BOM row 45: SMD push button ---
$ws.Range("B45").Value = "SMD Button"
$ws.Range("C45").Value = "SMD switch push button"
$ws.Range("D45").Value = "MJTPSMWBTR"
$ws.Range("E45").Value = "https://www.digikey.ca/en/products/detail/apem-inc/MJTPSMWBTR/1798006"
$ws.Range("F45").Value = 1.18
$ws.Range("G45").Value = 1

# Match the formatting used by the rest of the "total" column (fill style)
# before writing the formula, mirroring the existing shared-formula column.
$ws.Range("H44").Copy()
$ws.Range("H45").PasteSpecial(-4122)
$ws.Range("H45").Formula = "=G45*F45"

# Wire up the actual hyperlink, then give the new part's product-page cell
# the same "Hyperlink" look as its neighbours (order matters: Hyperlinks.Add
# applies its own style first, so re-apply "Hyperlink" afterwards).
$ws.Hyperlinks.Add($ws.Range("E45"), "https://www.digikey.ca/en/products/detail/apem-inc/MJTPSMWBTR/1798006")
$ws.Range("E45").Style = "Hyperlink"

# --- Update the saved view/selection ---
[void]$ws.Range("F3").Select()
